$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank columns at E..H (this shifts the existing "Created Date"
# column, and all its data, from E to I, leaving E..H empty for the new data)
$ws.Range("E1:H1").EntireColumn.Insert()

# Label the newly inserted columns with the "-Modified" header variants
$ws.Cells.Item(1, 5).Value = "User-Modified"
$ws.Cells.Item(1, 6).Value = "Is X Done?-Modified"
$ws.Cells.Item(1, 7).Value = "Is Y done?-Modified"
$ws.Cells.Item(1, 8).Value = "Is Z Done?-Modified"
